$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7475.1665
$ws.Range("I62").Value = 8000
$ws.Range("J62").Value = 7212.75
$ws.Range("K62").Value = 8000
$ws.Range("L62").Value = 7212.75
$ws.Range("M62").Value = -7376
$ws.Range("N62").Value = -8460.75
$ws.Range("H65").Value = 7475.1665
$ws.Range("I65").Value = 8000
$ws.Range("J65").Value = 7212.75
$ws.Range("K65").Value = 40000
$ws.Range("L65").Value = 36063.75
$ws.Range("M65").Value = -36880
$ws.Range("N65").Value = -42303.75
$ws.Range("H86").Value = 9457.076999999999
$ws.Range("I86").Value = 1766.6666
$ws.Range("J86").Value = 26760.5
$ws.Range("K86").Value = 1766.6666
$ws.Range("L86").Value = 26760.5
$ws.Range("M86").Value = -643.6666
$ws.Range("N86").Value = -29006.5
$ws.Range("H89").Value = 9457.076999999999
$ws.Range("I89").Value = 1766.6666
$ws.Range("J89").Value = 26760.5
$ws.Range("K89").Value = 8833.333000000001
$ws.Range("L89").Value = 133802.5
$ws.Range("M89").Value = -3217.333000000001
$ws.Range("N89").Value = -145034.5
$ws.Range("H98").Value = 707.8570999999999
$ws.Range("I98").Value = 685.38464
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 685.38464
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 812.61536
$ws.Range("N98").Value = -3996
$ws.Range("H116").Value = 16607392
$ws.Range("I116").Value = 40323936
$ws.Range("J116").Value = 5810.2
$ws.Range("K116").Value = 40323936
$ws.Range("L116").Value = 5810.2
$ws.Range("M116").Value = -40320494
$ws.Range("N116").Value = -12694.2
$ws.Range("H122").Value = 707.8570999999999
$ws.Range("I122").Value = 685.38464
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2056.15392
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 393.8460800000003
$ws.Range("N122").Value = -7900
$ws.Range("H129").Value = 271201.88
$ws.Range("I129").Value = 265.66666
$ws.Range("J129").Value = 295108
$ws.Range("K129").Value = 796.9999799999999
$ws.Range("L129").Value = 885324
$ws.Range("M129").Value = 4203.00002
$ws.Range("N129").Value = -895324
$ws.Range("H137").Value = 109977.836
$ws.Range("I137").Value = 155361.47
$ws.Range("J137").Value = 2707.4546
$ws.Range("K137").Value = 466084.41
$ws.Range("L137").Value = 8122.3638
$ws.Range("M137").Value = -463534.41
$ws.Range("N137").Value = -13222.3638
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 33335010
$ws.Range("I74").Value = 47619748
$ws.Range("K74").Value = 47619748
$ws.Range("M74").Value = -47618874
$ws.Range("H77").Value = 33335010
$ws.Range("I77").Value = 47619748
$ws.Range("K77").Value = 238098740
$ws.Range("M77").Value = -238094372
$ws.Range("H122").Value = 3502.9092
$ws.Range("I122").Value = 3559.2222
$ws.Range("K122").Value = 10677.6666
$ws.Range("M122").Value = -8227.6666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3925.379
$ws.Range("I31").Value = 1753.0333
$ws.Range("K31").Value = 1753.0333
$ws.Range("M31").Value = -1458.0333
$ws.Range("H34").Value = 3925.379
$ws.Range("I34").Value = 1753.0333
$ws.Range("K34").Value = 1753.0333
$ws.Range("M34").Value = -1551.0333
$ws.Range("H68").Value = 56495
$ws.Range("J68").Value = 56495
$ws.Range("L68").Value = 56495
$ws.Range("N68").Value = -57993
$ws.Range("H71").Value = 56495
$ws.Range("J71").Value = 56495
$ws.Range("L71").Value = 169485
$ws.Range("N71").Value = -176973
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H107").Value = 1173.4242
$ws.Range("I107").Value = 739.6111
$ws.Range("J107").Value = 1694
$ws.Range("K107").Value = 739.6111
$ws.Range("L107").Value = 1694
$ws.Range("M107").Value = 1180.3889
$ws.Range("N107").Value = -5534
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 845.8570999999999
$ws.Range("I113").Value = 671.25
$ws.Range("J113").Value = 953.3077
$ws.Range("K113").Value = 2013.75
$ws.Range("L113").Value = 2859.9231
$ws.Range("M113").Value = 156.25
$ws.Range("N113").Value = -7199.9231
$ws.Range("H131").Value = 661.77
$ws.Range("J131").Value = 693.5730600000001
$ws.Range("L131").Value = 2080.71918
$ws.Range("N131").Value = -12160.71918
$ws.Range("H137").Value = 16672631
$ws.Range("I137").Value = 1510
$ws.Range("J137").Value = 19614594
$ws.Range("K137").Value = 4530
$ws.Range("L137").Value = 58843782
$ws.Range("M137").Value = 570
$ws.Range("N137").Value = -58853982
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3369.276
$ws.Range("I80").Value = 2992.5386
$ws.Range("J80").Value = 3675.375
$ws.Range("K80").Value = 2992.5386
$ws.Range("L80").Value = 3675.375
$ws.Range("M80").Value = -1994.5386
$ws.Range("N80").Value = -5671.375
$ws.Range("H83").Value = 3369.276
$ws.Range("I83").Value = 2992.5386
$ws.Range("J83").Value = 3675.375
$ws.Range("K83").Value = 14962.693
$ws.Range("L83").Value = 18376.875
$ws.Range("M83").Value = -9970.692999999999
$ws.Range("N83").Value = -28360.875
$ws.Range("H102").Value = 8467.6
$ws.Range("I102").Value = 8081
$ws.Range("J102").Value = 10014
$ws.Range("K102").Value = 8081
$ws.Range("L102").Value = 10014
$ws.Range("M102").Value = -6459
$ws.Range("N102").Value = -13258
$ws.Range("H122").Value = 4461.5
$ws.Range("I122").Value = 4855.4287
$ws.Range("K122").Value = 14566.2861
$ws.Range("M122").Value = -12116.2861
$ws.Range("H126").Value = 3700
$ws.Range("I126").Value = 2756.25
$ws.Range("K126").Value = 8268.75
$ws.Range("M126").Value = -5798.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 90664
$ws.Range("I40").Value = 128314.664
$ws.Range("J40").Value = 5950
$ws.Range("K40").Value = 128314.664
$ws.Range("L40").Value = 5950
$ws.Range("M40").Value = -128178.664
$ws.Range("N40").Value = -6222
$ws.Range("H68").Value = 2833
$ws.Range("I68").Value = 2750
$ws.Range("J68").Value = 2999
$ws.Range("K68").Value = 2750
$ws.Range("L68").Value = 2999
$ws.Range("M68").Value = -2001
$ws.Range("N68").Value = -4497
$ws.Range("H71").Value = 2833
$ws.Range("I71").Value = 2750
$ws.Range("J71").Value = 2999
$ws.Range("K71").Value = 13750
$ws.Range("L71").Value = 14995
$ws.Range("M71").Value = -10006
$ws.Range("N71").Value = -22483
$ws.Range("H122").Value = 2805063.5
$ws.Range("I122").Value = 3924488.8
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 11773466.4
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -11771016.4
$ws.Range("N122").Value = -24400
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1846.8
$ws.Range("I126").Value = 1464.8334
$ws.Range("J126").Value = 2829
$ws.Range("K126").Value = 4394.5002
$ws.Range("L126").Value = 8487
$ws.Range("M126").Value = -1924.5002
$ws.Range("N126").Value = -13427
$ws.Range("H136").Value = 27780810
$ws.Range("I136").Value = 38463064
$ws.Range("K136").Value = 115389192
$ws.Range("M136").Value = -115386642
